$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 45

# Row 3
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 30

# Row 4
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 60

# Row 5
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 45

# Row 6
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 45

# Row 7 - C7 becomes the text "1.0"; H7 becomes empty text; I7 becomes numeric 45
$c7 = $ws.Range("C7")
$c7.Value = "'1.0"
$c7.Style = "Normal"

$h7 = $ws.Range("H7")
$h7.Value = "'"
$h7.Style = "Normal"

$ws.Range("I7").Value = 45

# Row 8 - C8 becomes the text "1.0"; H8 becomes empty text; I8 becomes numeric 45
$c8 = $ws.Range("C8")
$c8.Value = "'1.0"
$c8.Style = "Normal"

$h8 = $ws.Range("H8")
$h8.Value = "'"
$h8.Style = "Normal"

$ws.Range("I8").Value = 45

# Row 9
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 30

# Row 10
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 30

# Row 11
$ws.Range("I11").Value = 60
